$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: add a plain date (A8) and a fractional-day value (B8); the shared
# formula in C8 (=A8+B8) now evaluates against real data instead of blanks.
$ws.Range("A8").Value = 43171
$ws.Range("B8").Value = 0.846

# The old NOW() volatile formula in E8 is removed, leaving an empty (but
# still styled) cell.
$ws.Range("E8").ClearContents()

# Rows 9 and 10: add timestamp values in column A so the shared formula in
# column C picks them up.
$ws.Range("A9").Value = 43171.6525247569
$ws.Range("A10").Value = 43171.676977662

# Update the active selection to E9.
$ws.Range("E9").Select()
